$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7972
$ws.Range("B4").Value = 12680
$ws.Range("B5").Value = 20652
$ws.Range("B6").Value = 0.136765
